$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the Price (D) / Volume(1h) (E) columns store plain text that
# can look numeric (e.g. "1.000", "29.396.74"). Force each touched cell to
# Text format before writing so Excel does not auto-convert it to a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.396.74"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.848.89"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6277"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07633"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2908"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.72"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07742"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.032"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6784"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001062"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.24"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.149"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.408.87"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "226.70"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.32"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.482"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "157.89"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1380"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.411"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.67"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.389"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.467"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05587"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.125"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.056"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.835"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.6978"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.587"
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01803"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.229.79"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.719"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.404"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9029"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.63"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.89"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.168"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.980"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.676"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05702"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4631"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.14%  "
